$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lit review")

# Make "lit review" the active sheet (it was "lexical tasks" before).
$ws.Activate()

# Add a new row (23) documenting Warriner, Kuperman & Brysbaert (2013).
# (Cells are populated in the same order the original author typed them, so
# that new shared-string entries land at the same indices as the source file.)
$ws.Cells.Item(23, 1).Value2 = "Warriner, A. B., Kuperman, V., & Brysbaert, M. (2013). Norms of valence, arousal, and dominance for 13,915 English lemmas. Behavior Research Methods, 45(4), 1191–1207. https://doi.org/10.3758/s13428-012-0314-x"
$ws.Cells.Item(23, 7).Value2 = "It might be useful to analyze switch words based on valence strength (delta with neutral) as opposed to just binary neg2pos and pos2neg."
$ws.Cells.Item(23, 2).Value2 = "Affective ratings of ~14k lemmas to replace/extend ANEW, plus demographic and relational analyses."
$ws.Cells.Item(23, 3).Value2 = "Authors use mTurk to create a database of valence, arousal, and dominance ratings for English word lemmas in order to address the mismatch between the limited words available in the ANEW database and the megastudy analyses currently prevalent in psycholinguistics.  They also gather demographic information from the raters and provide information on primary correlations found."
$ws.Cells.Item(23, 4).Value2 = "On scales from low-to-high (unhappy/calm/controlled > happy/excited/in control), distributions of valence and dominance are negatively skewed: more words make people feel happy/in control (versus unhappy/controlled) whereas arousal is positively skewed: fewer words make people feel excited.  Overall, ratings of valence were relatively consistent across participants, but arousal and dominance had higher standard deviations, indicating more variability across participant responses."
$ws.Cells.Item(23, 5).Value2 = "Words high on the valence scale (that make people happy) have greater imageability, concreteness, familiarity, context availability, and body-object interaction; they are also higher in frequency and learned at earlier age.  Virtually all of these same words also make people feel in control."
$ws.Cells.Item(23, 6).Value2 = "Across all three dimensions, younger (v. older), lower (v. higher) education, and male (v. female) gave slightly higher ratings.  Female raters provided more extreme valence and dominance ratings for words on opposing ends of the frequency spectrum, leading a broader range of ratings for valence and dominance for female raters; the relationship between frequency and arousal for female raters was weak."
$ws.Cells.Item(23, 8).Value2 = "`"The number of words covered by the ANEW norms appeared sufficient for use in small-scale factorial experiments. In these experiments, a limited number of stimuli would be selected that varied on one dimension (e.g., valence) and were matched on other variables (e.g., arousal, word frequency, and word length). However, the number of words in this set is prohibitively small for the large-scale megastudies that are currently emerging in psycholinguistics.`"
`"To sum up, in terms of the variability of ratings, valence and dominance pattern together and are best considered in terms of their magnitude (how strong is the feeling) rather than their polarity (sad vs. happy, or controlled by vs. in control); polarity, however, determines variability in the arousal ratings.`"
`"The fact that extreme values of valence and dominance are more arousing point again at the utility of considering valence/dominance strength (i.e., how different a word is from neutral) rather than polarity as the explanatory variable.`""

# Match formatting of the rest of the table (wrapped text, thin borders, row height).
$rng = $ws.Range("A23:H23")
$rng.WrapText = $true
$rng.Borders.LineStyle = 1
$ws.Rows.Item(23).RowHeight = 380

# Leave the selection on the last cell touched, as in the saved file.
$ws.Range("G23").Select() | Out-Null
